$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.007.27"
$ws.Range("E2").Value = "  -2.83%  "
$ws.Range("D3").Value = "2.339.47"
$ws.Range("E3").Value = "  -3.91%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.48"
$ws.Range("E5").Value = "  -1.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "85.29"
$ws.Range("E6").Value = "  -4.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.529"
$ws.Range("E7").Value = "  -2.20%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -2.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0811"
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "29.98"
$ws.Range("E11").Value = "  -6.71%  "
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").Value = "2.699.36"
$ws.Range("E13").Value = "  -3.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.43"
$ws.Range("E14").Value = "  -4.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.75"
$ws.Range("E15").Value = "  -5.39%  "
$ws.Range("D16").Value = "2.368.52"
$ws.Range("E16").Value = "  -2.60%  "
$ws.Range("E17").Value = "  -1.94%  "
$ws.Range("D18").Value = "39.990.98"
$ws.Range("E18").Value = "  -2.73%  "
$ws.Range("D19").Value = "0.0₃0903"
$ws.Range("E19").Value = "  -2.22%  "
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.03"
$ws.Range("E21").Value = "  -5.57%  "
$ws.Range("E22").Value = "  -2.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.18"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.56"
$ws.Range("E24").Value = "  -4.62%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  -3.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.33"
$ws.Range("E27").Value = "  -2.89%  "
$ws.Range("E28").Value = "  -4.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.30"
$ws.Range("E29").Value = "  -2.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.75"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "152.97"
$ws.Range("E31").Value = "  -2.22%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.11"
$ws.Range("E33").Value = "  -2.97%  "
$ws.Range("E34").Value = "  -2.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0718"
$ws.Range("E35").Value = "  -3.50%  "
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("E37").Value = "  -3.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0986"
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("E39").Value = "  -6.07%  "
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.88"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("D42").Value = "1.954.07"
$ws.Range("E42").Value = "  -1.67%  "
$ws.Range("E43").Value = "  -4.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.77"
$ws.Range("E44").Value = "  -3.35%  "
$ws.Range("E45").Value = "  -4.23%  "
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("E47").Value = "  -6.00%  "
$ws.Range("D48").Value = "2.557.69"
$ws.Range("E48").Value = "  -4.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "92.76"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.58"
$ws.Range("E50").Value = "  -3.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.82"
$ws.Range("E51").Value = "  -1.26%  "
